$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.301.34'
$ws.Range("E2").Value = '  -1.84%  '
$ws.Range("D3").Value = '2.446.20'
$ws.Range("E3").Value = '  -1.51%  '
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.19'
$ws.Range("E5").Value = '  -2.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.48'
$ws.Range("E6").Value = '  -1.74%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.500'
$ws.Range("E8").Value = '  -2.15%  '
$ws.Range("D9").Value = '2.445.74'
$ws.Range("E9").Value = '  -1.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.148'
$ws.Range("E10").Value = '  -6.18%  '
$ws.Range("E11").Value = '  -1.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.334'
$ws.Range("E12").Value = '  -5.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.75'
$ws.Range("E13").Value = '  -3.24%  '
$ws.Range("D14").Value = '2.898.97'
$ws.Range("E14").Value = '  -1.56%  '
$ws.Range("D15").Value = '68.240.53'
$ws.Range("E15").Value = '  -1.82%  '
$ws.Range("E16").Value = '  -4.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '23.04'
$ws.Range("E17").Value = '  -4.93%  '
$ws.Range("D18").Value = '2.459.38'
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.80'
$ws.Range("E19").Value = '  -3.34%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '338.66'
$ws.Range("E20").Value = '  -1.85%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.08'
$ws.Range("E21").Value = '  -3.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.73'
$ws.Range("E22").Value = '  -3.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.83'
$ws.Range("E24").Value = '  -4.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.11'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.63'
$ws.Range("E26").Value = '  -6.33%  '
$ws.Range("D27").Value = '2.574.81'
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.02'
$ws.Range("E29").Value = '  -7.12%  '
$ws.Range("D30").Value = '0.0₃0818'
$ws.Range("E30").Value = '  -6.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.09'
$ws.Range("E31").Value = '  -9.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '422.87'
$ws.Range("E33").Value = '  -4.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.14'
$ws.Range("E34").Value = '  -4.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.63'
$ws.Range("E35").Value = '  -4.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '157.08'
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.99'
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.108'
$ws.Range("E39").Value = '  -4.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.69'
$ws.Range("E40").Value = '  -2.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.301'
$ws.Range("E41").Value = '  -4.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.38'
$ws.Range("E42").Value = '  -4.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.48'
$ws.Range("E43").Value = '  -6.25%  '
$ws.Range("E44").Value = '  +0.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '133.33'
$ws.Range("E45").Value = '  -4.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.02'
$ws.Range("E46").Value = '  -6.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.31'
$ws.Range("E47").Value = '  -3.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0711'
$ws.Range("E48").Value = '  -2.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.475'
$ws.Range("E49").Value = '  -7.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.557'
$ws.Range("E50").Value = '  -2.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0903'
$ws.Range("E51").Value = '  -1.67%  '
